$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: 1K_PFF_HTRA1_Hoechst_20x_04
$ws.Range("B5").Value = 6
$ws.Range("E5").Value = 2099
$ws.Range("F5").Value = 0

# Row 8: 1K_PFF_HTRA1_Hoechst_20x_07
$ws.Range("B8").Value = 16
$ws.Range("E8").Value = 1162
$ws.Range("F8").Value = 82

# Row 12: 1K_PFF_noHTRA1_Hoechst_20x_01
$ws.Range("D12").Value = 0
$ws.Range("F12").Value = 0

# Row 13: 1K_PFF_noHTRA1_Hoechst_20x_02
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 0

# Row 14: 1K_PFF_noHTRA1_Hoechst_20x_03
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0

# Row 15: 1K_PFF_noHTRA1_Hoechst_20x_04
$ws.Range("D15").Value = 0
$ws.Range("F15").Value = 0

# Row 16: 1K_PFF_noHTRA1_Hoechst_20x_05
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 0

# Row 17: 1K_PFF_noHTRA1_Hoechst_20x_06
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0

# Row 18: 1K_PFF_noHTRA1_Hoechst_20x_07
$ws.Range("D18").Value = 0
$ws.Range("F18").Value = 0

# Row 19: 1K_PFF_noHTRA1_Hoechst_20x_08
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0

# Row 20: 1K_PFF_noHTRA1_Hoechst_20x_09
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = 0

# Row 21: 1K_PFF_noHTRA1_Hoechst_20x_10
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = 0

# Row 23: WT_PFF_HTRA1_Hoechst_20x_02
$ws.Range("B23").Value = 19
$ws.Range("E23").Value = 2670
$ws.Range("F23").Value = 431

# Row 32: WT_PFF_noHTRA1_Hoechst_20x_01
$ws.Range("D32").Value = 0
$ws.Range("F32").Value = 0

# Row 33: WT_PFF_noHTRA1_Hoechst_20x_02
$ws.Range("D33").Value = 0
$ws.Range("F33").Value = 0

# Row 34: WT_PFF_noHTRA1_Hoechst_20x_03
$ws.Range("D34").Value = 0
$ws.Range("F34").Value = 0

# Row 35: WT_PFF_noHTRA1_Hoechst_20x_04
$ws.Range("D35").Value = 0
$ws.Range("F35").Value = 0

# Row 36: WT_PFF_noHTRA1_Hoechst_20x_05
$ws.Range("D36").Value = 0
$ws.Range("F36").Value = 0

# Row 37: WT_PFF_noHTRA1_Hoechst_20x_06
$ws.Range("D37").Value = 0
$ws.Range("F37").Value = 0

# Row 38: WT_PFF_noHTRA1_Hoechst_20x_07
$ws.Range("D38").Value = 0
$ws.Range("F38").Value = 0

# Row 39: WT_PFF_noHTRA1_Hoechst_20x_08
$ws.Range("D39").Value = 0
$ws.Range("F39").Value = 0

# Row 40: WT_PFF_noHTRA1_Hoechst_20x_09
$ws.Range("D40").Value = 0
$ws.Range("F40").Value = 0

# Row 41: WT_PFF_noHTRA1_Hoechst_20x_10
$ws.Range("D41").Value = 0
$ws.Range("F41").Value = 0
